# Updates the "NEW" sheet of the interactive map workbook:
#  - inserts a new incident row (Caso 6173 / ARMENIA 2321) at row 29,
#    shifting the existing rows 29-46 down to 30-47
#  - appends a new incident row (Caso 6171 / CABELLO 3486) as the new last row 48

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values look numeric/date-like and would otherwise be
# auto-converted by Excel (we need plain text, like the rest of the sheet).
$textCols = @(1, 2, 4, 5, 9)

function Set-RowValues($rowIndex, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 1
        $cell = $ws.Cells.Item($rowIndex, $col)
        if ($textCols -contains $col) {
            # Force plain-text storage, then strip the temporary formatting
            # so the saved cell carries no style index (matches sheet norm).
            $cell.NumberFormat = "@"
            $cell.Value = $values[$i]
            $cell.ClearFormats()
        } else {
            $cell.Value = $values[$i]
        }
    }
}

# --- Insert the new row 29 (push old rows 29-46 down to 30-47) ---
$ws.Rows.Item(29).Insert()

$row29 = @("6173", "4/29/2025", "ARMENIA 2321", "14", "805507398", "NEW", "Pendiente", "Picada", "1", "Cambio", "Sin equipos", "Pasante")
Set-RowValues 29 $row29
$ws.Cells.Item(29, 13).Value = -58.420549
$ws.Cells.Item(29, 14).Value = -34.585103

# --- Append the new last row 48 ---
$row48 = @("6171", "6/18/2025", "CABELLO 3486", "14", "807658640", "NEW", "Pendiente", "Columna inclinada evaluar con inspector un corrimiento", "1", "Cambio", "Sin equipos", "Pasante")
Set-RowValues 48 $row48
$ws.Cells.Item(48, 13).Value = -58.409579
$ws.Cells.Item(48, 14).Value = -34.581134
